# The workbook lists master "device_type" rows grouped in triplets
# (eng / ara / fra) per device code: FRS, IRS, DKS (Desktop), CMR, SCN, PRT.
# This edit removes the "DKS" (Desktop Computer) triplet of rows, which
# currently occupies rows 8-10 of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three "Desktop" rows (code = DKS): eng/ara/fra entries.
$ws.Rows("8:10").Delete() | Out-Null

# Configure the page for printing (paperSize 9 = A4, orientation 1 = portrait).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Leave the active selection on E10, matching where the cursor ended up
# after the row deletion.
$ws.Range("E10").Select() | Out-Null
